# Change the color of the "dwc:Taxon" class (and its surrounding oval) in the
# ACS diagram from blue (00B0F0) to red (FF0000).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Oval 12 is the ellipse drawn around the "dwc:Taxon" label.
$oval = $s.Shapes.Item(11)
$oval.Line.ForeColor.RGB = 255

# TextBox 103 holds the "dwc:Taxon" label itself (originally split across two
# runs - "dwc" and ":Taxon" - both colored blue). Re-enter the text as a
# single run and recolor/restyle it red. Turn off the shape's "shrink/grow
# to fit text" behavior while we do this so its box doesn't get resized by
# the momentary text clear.
$tb = $s.Shapes.Item(45)
$tb.TextFrame.AutoSize = 0
$tb.TextFrame.TextRange.Font.Color.RGB = 255
$tb.TextFrame.TextRange.Text = ""
$tb.TextFrame.TextRange.Text = "dwc:Taxon"
$tb.TextFrame.TextRange.Font.Size = 36
$tb.TextFrame.TextRange.Font.Bold = $true
$tb.TextFrame.TextRange.Font.Color.RGB = 255
$tb.TextFrame.AutoSize = 1
